$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (5-55): refreshed stdev values ---
$ws.Range("B5").Value = 5.0561127662658691
$ws.Range("C5").Value = 5.0083060264587402
$ws.Range("E5").Value = 5.0762591361999512
$ws.Range("C6").Value = 5.180595874786377
$ws.Range("E6").Value = 5.0988507270812988
$ws.Range("C7").Value = 5.2569370269775391
$ws.Range("E7").Value = 5.1142368316650391
$ws.Range("C8").Value = 5.2174882888793945
$ws.Range("E8").Value = 5.0969171524047852
$ws.Range("B9").Value = 4.8313360214233398
$ws.Range("C9").Value = 4.7179679870605469
$ws.Range("E9").Value = 5.1106996536254883
$ws.Range("C10").Value = 5.2118101119995117
$ws.Range("E10").Value = 5.1036257743835449
$ws.Range("C11").Value = 5.2065525054931641
$ws.Range("E11").Value = 5.0792875289916992
$ws.Range("C12").Value = 4.9756793975830078
$ws.Range("E12").Value = 5.0889816284179688
$ws.Range("C13").Value = 5.2209606170654297
$ws.Range("E13").Value = 5.0488548278808594
$ws.Range("C14").Value = 4.944638729095459
$ws.Range("E14").Value = 5.0137906074523926
$ws.Range("C15").Value = 4.9615521430969238
$ws.Range("E15").Value = 5.0443763732910156
$ws.Range("C16").Value = 5.3441829681396484
$ws.Range("E16").Value = 5.0197534561157227
$ws.Range("C17").Value = 4.8563470840454102
$ws.Range("E17").Value = 4.994544506072998
$ws.Range("C18").Value = 4.4023919105529785
$ws.Range("E18").Value = 4.9905996322631836
$ws.Range("C19").Value = 5.4870815277099609
$ws.Range("E19").Value = 4.9392385482788086
$ws.Range("C20").Value = 4.9849457740783691
$ws.Range("E20").Value = 4.9605741500854492
$ws.Range("C21").Value = 4.7488002777099609
$ws.Range("E21").Value = 4.9281167984008789
$ws.Range("C22").Value = 5.185457706451416
$ws.Range("D22").Value = 5.0657110214233398
$ws.Range("E22").Value = 5.0045766830444336
$ws.Range("C23").Value = 4.4823884963989258
$ws.Range("E23").Value = 5.0837841033935547
$ws.Range("C24").Value = 5.1535735130310059
$ws.Range("E24").Value = 5.1259121894836426
$ws.Range("C25").Value = 5.0520658493041992
$ws.Range("E25").Value = 5.2856388092041016
$ws.Range("C26").Value = 5.5444831848144531
$ws.Range("E26").Value = 5.4715142250061035
$ws.Range("C27").Value = 5.1152591705322266
$ws.Range("D27").Value = 5.630683422088623
$ws.Range("E27").Value = 5.5879006385803223
$ws.Range("C28").Value = 5.8662357330322266
$ws.Range("E28").Value = 5.8015623092651367
$ws.Range("C29").Value = 6.4224839210510254
$ws.Range("E29").Value = 5.9782395362854004
$ws.Range("C30").Value = 6.4216823577880859
$ws.Range("E30").Value = 6.2029051780700684
$ws.Range("C31").Value = 6.2329339981079102
$ws.Range("E31").Value = 6.3807382583618164
$ws.Range("C32").Value = 6.4053440093994141
$ws.Range("E32").Value = 6.557830810546875
$ws.Range("C33").Value = 6.7436680793762207
$ws.Range("E33").Value = 6.6444807052612305
$ws.Range("C34").Value = 7.0740551948547363
$ws.Range("E34").Value = 6.689420223236084
$ws.Range("C35").Value = 7.1449832916259766
$ws.Range("E35").Value = 6.7020959854125977
$ws.Range("C36").Value = 6.7090897560119629
$ws.Range("E36").Value = 6.7660746574401855
$ws.Range("C37").Value = 6.6460857391357422
$ws.Range("E37").Value = 6.7728581428527832
$ws.Range("C38").Value = 6.8269381523132324
$ws.Range("E38").Value = 6.7536067962646484
$ws.Range("C39").Value = 6.5357670783996582
$ws.Range("E39").Value = 6.6841130256652832
$ws.Range("C40").Value = 6.8087425231933594
$ws.Range("D40").Value = 6.6376562118530273
$ws.Range("E40").Value = 6.587486743927002
$ws.Range("C41").Value = 6.4663949012756348
$ws.Range("E41").Value = 6.5504288673400879
$ws.Range("C42").Value = 6.5704026222229004
$ws.Range("E42").Value = 6.463472843170166
$ws.Range("C43").Value = 6.4486150741577148
$ws.Range("E43").Value = 6.3613138198852539
$ws.Range("C44").Value = 6.2753458023071289
$ws.Range("E44").Value = 6.2453608512878418
$ws.Range("C45").Value = 6.3755683898925781
$ws.Range("E45").Value = 6.0911006927490234
$ws.Range("C46").Value = 5.8634829521179199
$ws.Range("E46").Value = 6.001129150390625
$ws.Range("C47").Value = 5.9075040817260742
$ws.Range("E47").Value = 5.9087715148925781
$ws.Range("C48").Value = 5.492192268371582
$ws.Range("E48").Value = 5.8488078117370605
$ws.Range("C49").Value = 5.4204001426696777
$ws.Range("E49").Value = 5.8128242492675781
$ws.Range("B50").Value = 5.6874308586120605
$ws.Range("C50").Value = 5.6566510200500488
$ws.Range("E50").Value = 5.7273802757263184
$ws.Range("C51").Value = 5.7391834259033203
$ws.Range("E51").Value = 5.7291460037231445
$ws.Range("C52").Value = 5.9089431762695312
$ws.Range("D52").Value = 5.755403995513916
$ws.Range("E52").Value = 5.6727848052978516
$ws.Range("C53").Value = 5.9514932632446289
$ws.Range("D53").Value = 5.7368502616882324
$ws.Range("E53").Value = 5.6481046676635742
$ws.Range("C54").Value = 5.6065735816955566
$ws.Range("D54").Value = 5.7360215187072754
$ws.Range("E54").Value = 5.6395883560180664
$ws.Range("C55").Value = 5.8793721199035645
$ws.Range("D55").Value = 5.7184667587280273
$ws.Range("E55").Value = 5.614771842956543

# --- Fill previously-empty rows 56-57 (already styled) ---
$ws.Range("B56").Value = 5.5315570831298828
$ws.Range("C56").Value = 5.4002552032470703
$ws.Range("D56").Value = 5.6665811538696289
$ws.Range("E56").Value = 5.5613002777099609
$ws.Range("B57").Value = 5.3835568428039551
$ws.Range("C57").Value = 5.2700705528259277
$ws.Range("D57").Value = 5.5828619003295898
$ws.Range("E57").Value = 5.4797205924987793

# --- Add new rows 58-64: copy formatting from row 55, then set values ---
$ws.Range("A55:E55").Copy()
$ws.Range("A58:E64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A58").Value = 45536
$ws.Range("B58").Value = 5.4593696594238281
$ws.Range("C58").Value = 5.3437538146972656
$ws.Range("D58").Value = 5.5301661491394043
$ws.Range("E58").Value = 5.4284706115722656
$ws.Range("A59").Value = 45566
$ws.Range("B59").Value = 5.5294394493103027
$ws.Range("C59").Value = 5.4333004951477051
$ws.Range("D59").Value = 5.5471291542053223
$ws.Range("E59").Value = 5.4447407722473145
$ws.Range("A60").Value = 45597
$ws.Range("B60").Value = 5.3712029457092285
$ws.Range("C60").Value = 5.2579398155212402
$ws.Range("D60").Value = 5.5099964141845703
$ws.Range("E60").Value = 5.4059686660766602
$ws.Range("A61").Value = 45627
$ws.Range("B61").Value = 5.2135939598083496
$ws.Range("C61").Value = 5.1747245788574219
$ws.Range("D61").Value = 5.5073013305664062
$ws.Range("E61").Value = 5.4066824913024902
$ws.Range("A62").Value = 45658
$ws.Range("B62").Value = 5.5685210227966309
$ws.Range("C62").Value = 5.4902467727661133
$ws.Range("D62").Value = 5.5249791145324707
$ws.Range("E62").Value = 5.4261984825134277
$ws.Range("A63").Value = 45689
$ws.Range("B63").Value = 5.8855748176574707
$ws.Range("C63").Value = 5.7530045509338379
$ws.Range("D63").Value = 5.5359139442443848
$ws.Range("E63").Value = 5.4399394989013672
$ws.Range("A64").Value = 45717
$ws.Range("B64").Value = 5.6471514701843262
$ws.Range("C64").Value = 5.5304203033447266
$ws.Range("D64").Value = 5.5372090339660645
$ws.Range("E64").Value = 5.4412670135498047
